$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Choice commitment" -> "Choice "
$ws.Range("A43").Value = 'Choice '

# Row labels: "Forced commitment " -> "Mandatory structured"
$ws.Range("A5").Value = 'Mandatory structured'
$ws.Range("A14").Value = 'Mandatory structured'
$ws.Range("A23").Value = 'Mandatory structured'
$ws.Range("A32").Value = 'Mandatory structured'
$ws.Range("A41").Value = 'Mandatory structured'

# Panel header rows: reword "Forced Commitment" -> "Mandatory structured"
$ws.Range("B3").Value = 'Panel A : $\quad$ Control  = 0           $\quad\quad$                 Mandatory structured  = 0'
$ws.Range("B12").Value = 'Panel B : $\quad$ Control  = 0         $\quad\quad$                    Mandatory structured = 1'
$ws.Range("B21").Value = 'Panel C : $\quad$ Control  = 1        $\quad\quad$                     Mandatory structured = 0'
$ws.Range("B30").Value = 'Panel D : $\quad$ Control  = 1       $\quad\quad$                      Mandatory structured = 1'

# Restore selection / view state
$ws.Range("G48").Select()
